# PlayerPerformance_4692.xlsx edit:
#  1. Insert a new "Player Info" sheet at the front with ID/NAME/BATTING_HAND/BOWL_STYLE data.
#  2. On "ODI Batting" and "ODI Bowling", rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#     full howstat.com scorecard URL values with just the bare numeric match code.

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# 1. New "Player Info" sheet (added first, so it lands before the existing
#    two sheets - matching sheetId order 1/2/3 in the target workbook).
# -------------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add()
$infoSheet.Name = "Player Info"

# Header row, styled like the header rows on the other sheets (bold, centered,
# thin box border).
$header = $infoSheet.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

# Data row - force text ("'" prefix) so "4692" is stored as a string, like
# every other "numeric" cell value in this workbook.
$infoSheet.Range("A2").Value = "'4692"
$infoSheet.Range("B2").Value = "Liam Stephen Livingstone"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Right Arm Leg Break"

# -------------------------------------------------------------------------
# 2. ODI Batting sheet: MATCH_CARD_LINK column -> MATCH_CODE, values become
#    just the bare match-code number instead of the full scorecard URL.
# -------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2  = "4456"
    3  = "4457"
    4  = "4469"
    5  = "4598"
    6  = "4599"
    7  = "4602"
    8  = "4609"
    9  = "4613"
    10 = "4618"
    11 = "4619"
    12 = "4620"
    13 = "4622"
}
foreach ($row in $battingCodes.Keys) {
    $battingSheet.Range("D$row").Value = "'" + $battingCodes[$row]
}

# -------------------------------------------------------------------------
# 3. ODI Bowling sheet: same MATCH_CARD_LINK -> MATCH_CODE treatment, but the
#    link lives in column B here.
# -------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "4457"
    3 = "4599"
    4 = "4602"
    5 = "4613"
    6 = "4618"
    7 = "4619"
    8 = "4622"
}
foreach ($row in $bowlingCodes.Keys) {
    $bowlingSheet.Range("B$row").Value = "'" + $bowlingCodes[$row]
}
